# Add a new "28-ago" column (BG) to the right of the existing "25-ago" column (BF)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new date column
$ws.Range("BG1").Value = "28-ago"

# Data values for the new column, one per data row (rows 2-18)
$values = @(
    0,
    14.097377768918269,
    20.720229759568294,
    24.122039928300879,
    0,
    9.9124666536503501,
    23.256071045566348,
    9.9692034840391344,
    0.84658676151751278,
    15.034121883031451,
    0,
    14.465597790691803,
    0,
    0,
    39.157878186329093,
    0,
    0
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 59).Value = $values[$i]
}

# Update the active selection to match the new edit location
$ws.Range("BI7").Select()
